$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "sr.no." column (A) for existing rows 4-8, which was previously blank
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# Add two new rows of data: build() and perform()
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "build()"
$ws.Range("C9").Value = "method is used to compile all the listed actions into a single step"
$ws.Range("D9").Value = "no argument"
$ws.Range("E9").Value = "actions class object"

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "perform()"
$ws.Range("C10").Value = "used to perform the action "
$ws.Range("D10").Value = "no argument"
$ws.Range("E10").Value = "actions class object"

# Update the selection to match the final state seen in the diff
$ws.Range("E9:E10").Select()

$wb.Save()
